# Apply the edit described by the diff:
#  1. Rename the worksheet/sheet-name entry from "Gamma1F-HW25.xpc" to "Gamma1F"
#  2. Append a new row 16 with averaged-intensity data for the
#     "HexGrid-60degTilt5degRes" quadrature scheme (reusing the existing
#     shared string used by row 15, index 13 -> "HexGrid-60degTilt5degRes")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the sheet/tab
$ws.Name = "Gamma1F"

# 2) Append the new data row (row 16)
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 1.089175334079616
$ws.Range("D16").Value = 0.7561887107447419
$ws.Range("E16").Value = 1.036456733616979
$ws.Range("F16").Value = 1.089175334079616
$ws.Range("G16").Value = 0.8664604101021188
$ws.Range("H16").Value = 1.102811925910289
$ws.Range("I16").Value = 1.05363254966033
$ws.Range("J16").Value = 0.7561887107447419
$ws.Range("K16").Value = 0.8963227221808603
$ws.Range("L16").Value = 0.9927490281302379
$ws.Range("M16").Value = 0.9841209440190125

# Match the style of column A in other data rows (bordered header-ish style s="1")
# by copying formatting from the cell directly above (A15) instead of
# re-declaring individual format properties (avoids minting a redundant xf).
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
